# Anexo 12.1: replace the placeholder "{nombreEmpresa}" with "{empresa}"
# in the "De:" paragraph (e.g. "... Tutor de la Empresa {nombreEmpresa}"
# becomes "... Tutor de la Empresa {empresa}").

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "nombreEmpresa",   # FindText
    $true,             # MatchCase
    $false,            # MatchWholeWord
    $false,            # MatchWildcards
    $false,            # MatchSoundsLike
    $false,            # MatchAllWordForms
    $true,             # Forward
    1,                 # Wrap (wdFindContinue)
    $false,            # Format
    "empresa",         # ReplaceWith
    2                  # Replace (wdReplaceAll)
)
